$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column D to give more room for the new explanation text
# (target stored width ~= 47.21875; Excel's internal pixel rounding means this
# COM input value lands closest to it)
$ws.Columns.Item(4).ColumnWidth = 46.33

# Add an explanation for each "Recommended Action" row in the new column D
$ws.Range("D2").Value = "SME-QA (Subject Matter Expert Qualification Assessments) is a program in which subject matter experts partner closely with HR representatives to hire specialized talent. SMEs assess talent levels earlier in the hiring process, which leads to an applicant pool that more closely fits the needs of the SMEs and speeds up the hiring process. This program is an effective solution for APG teams looking to quickly hire specialized talent to move their goals forward."
$ws.Range("D3").Value = "IPM allows agencies to receive temporary personnel assignments. This program is specifically focused on short-term engagements of non-Federal workers in the Federal space. Assignments may come to or from state and local governments, institutions of higher education, Indian tribal governments and other eligible organizations. IPM is useful to agencies looking for team members to bring innovation and perspective from outside the Federal government and provide a valuable experience for a non-Federal worker."
$ws.Range("D4").Value = "18F is a Federal digital consultancy group housed within the GSA. 18F works with agencies to modernize and improve efficiency within their technical solutions while prioritizing the goals of the partnering agency and keeping long-term use in mind.  A partnership with 18F is ideal for groups looking for an experienced, innovative group of technologists within government who are exceptional at applying technological solutions to progress Federal goals."
$ws.Range("D5").Value = "TMF is a funding vehicle that seeks out applicants looking to innovate and modernize technology within government. Applicants awarded funding will receive funding on the completion of major project milestones. The TMF is an investment in an agency to develop a solution that will help provide better services to stakeholders and further the mission of the agency. This solution is a good option for agencies who have vision for their technical solutions, but lack the funding and/or technical experts to implement on that vision."

# Match the wrap-text styling used by the rest of the data cells
$ws.Range("D2:D5").WrapText = $true

# Grow the rows so the much longer wrapped text fits
$ws.Rows.Item(2).RowHeight = 216
$ws.Rows.Item(3).RowHeight = 244.8
$ws.Rows.Item(4).RowHeight = 129.6
$ws.Rows.Item(5).RowHeight = 144

# Update the view: scroll down and select D6, matching the final cursor position
$ws.Range("D6").Select()
$excel.ActiveWindow.ScrollRow = 4
